# Add the 4-Aug (serial 44777) measurement day as a brand-new first day of
# data, ahead of the existing 6-Aug (serial 44779) day. The six pots x three
# replicate rows each (18 rows total) are inserted right after the header
# row, which pushes every pre-existing row down by 18 (old row N -> new row
# N+18). Only "weight" (column G) is known for this new day; "mean_moisture"
# (C) and "sd" (D) are left blank, matching the pattern already used for the
# other un-processed days at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New day's rows: pot, soil_type order, irrigation (F) and weight (G).
$newRows = @(
    @(1, 1, 44777, 60, 755),
    @(1, 2, 44777, 90, 750),
    @(1, 3, 44777, 120, 740),
    @(1, 4, 44777, 60, 800),
    @(1, 5, 44777, 90, 750),
    @(1, 6, 44777, 120, 850),
    @(2, 1, 44777, 60, 610),
    @(2, 2, 44777, 90, 630),
    @(2, 3, 44777, 120, 570),
    @(2, 4, 44777, 60, 600),
    @(2, 5, 44777, 90, 625),
    @(2, 6, 44777, 120, 605),
    @(3, 1, 44777, 60, 552),
    @(3, 2, 44777, 90, 590),
    @(3, 3, 44777, 120, 545),
    @(3, 4, 44777, 60, 590),
    @(3, 5, 44777, 90, 550),
    @(3, 6, 44777, 120, 570)
)

$insertCount = $newRows.Count
$firstRow = 2
$lastRow = $firstRow + $insertCount - 1

# Push the existing data (rows 2..145) down by 18 rows.
$ws.Rows("$firstRow`:$lastRow").Insert()

# The freshly inserted rows inherited the header row's formatting (style
# "1" on every column). Re-apply the plain data-row formatting by copying
# it down from the row immediately below the inserted block (row 20, which
# now holds what used to be row 2 and already has the right per-column
# styles: default on A/B/F/G, style "1" on C/D, the date style on E).
$templateRow = $lastRow + 1
$ws.Range("A$templateRow`:G$templateRow").Copy() | Out-Null
$ws.Range("A$firstRow`:G$lastRow").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Write the new day's values into the inserted rows.
for ($i = 0; $i -lt $insertCount; $i++) {
    $r = $firstRow + $i
    $vals = $newRows[$i]
    $ws.Cells.Item($r, 1).Value2 = $vals[0]  # A: soil_type
    $ws.Cells.Item($r, 2).Value2 = $vals[1]  # B: pot
    # C (mean_moisture) and D (sd) stay blank for this day.
    $ws.Cells.Item($r, 5).Value2 = $vals[2]  # E: date
    $ws.Cells.Item($r, 6).Value2 = $vals[3]  # F: irrigation
    $ws.Cells.Item($r, 7).Value2 = $vals[4]  # G: weight
}

# Tidy up the view: drop the scrolled-down top-left cell and move the
# selection to reflect the newly added data.
$ws.Activate()
$ws.Range("G20").Select() | Out-Null
$wv = $excel.ActiveWindow
$wv.ScrollRow = 1
$wv.ScrollColumn = 1
